$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate the "GF #" column (J) with the closing/GF file numbers. Rows
# that did not close in the title company's system are marked with the
# text "did not close" instead of a GF number.
$gfValues = @{
    2  = 20191011
    3  = 20191062
    4  = "did not close"
    5  = "did not close"
    6  = "did not close"
    7  = "did not close"
    8  = "did not close"
    9  = "did not close"
    10 = 20182625
    11 = "did not close"
}

foreach ($row in $gfValues.Keys) {
    $ws.Cells.Item($row, 10).Value = $gfValues[$row]
}

$wb.Save()

Write-Host "Done checking which files closed!"
